$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before C: shifts old "cpu total time [sec]" (C) to D,
# and frees up C for the new "total time [sec]" column.
$ws.Range("C1").EntireColumn.Insert()

# Remove the now-obsolete last two rows (30 and 31).
$ws.Range("A30:D31").EntireRow.Delete()

# Header row
$ws.Range("A1").Value = "slopy%"
$ws.Range("B1").Value = "area [km^2]"
$ws.Range("C1").Value = "total time [sec]"
$ws.Range("D1").Value = "cpu total time [sec]"

# Row 2
$ws.Range("A2").Value = 39.06
$ws.Range("B2").Value = 0.04000000000000001
$ws.Range("C2").Value = 2.811877250671387
$ws.Range("D2").ClearContents()

# Row 3
$ws.Range("A3").Value = 39.06
$ws.Range("B3").Value = 0.04000000000000001
$ws.Range("C3").Value = 2.833401441574097
$ws.Range("D3").ClearContents()

# Row 4
$ws.Range("A4").Value = 39.06
$ws.Range("B4").Value = 0.04000000000000001
$ws.Range("C4").Value = 6.14512300491333
$ws.Range("D4").ClearContents()

# Row 5
$ws.Range("A5").Value = 39.78
$ws.Range("B5").Value = 0.16
$ws.Range("C5").Value = 10.57961821556091
$ws.Range("D5").Value = 9.328125

# Row 6
$ws.Range("A6").Value = 39.78
$ws.Range("B6").Value = 0.16
$ws.Range("C6").Value = 13.87360000610352
$ws.Range("D6").Value = 12.53125

# Row 7
$ws.Range("A7").Value = 39.06
$ws.Range("B7").Value = 0.04000000000000001
$ws.Range("C7").Value = 3.679166555404663
$ws.Range("D7").Value = 2.453125

# Row 8
$ws.Range("A8").Value = 39.06
$ws.Range("B8").Value = 0.04000000000000001
$ws.Range("C8").Value = 3.316033124923706
$ws.Range("D8").Value = 2.5625

# Row 9
$ws.Range("A9").Value = 39.06
$ws.Range("B9").Value = 0.04000000000000001
$ws.Range("C9").Value = 3.487463235855103
$ws.Range("D9").Value = 2.953125

# Row 10
$ws.Range("A10").Value = 39.06
$ws.Range("B10").Value = 0.04000000000000001
$ws.Range("C10").Value = 3.183549880981445
$ws.Range("D10").Value = 2.546875

# Row 11
$ws.Range("A11").Value = 39.78
$ws.Range("B11").Value = 0.16
$ws.Range("C11").Value = 10.6967921257019
$ws.Range("D11").Value = 9.96875

# Row 12
$ws.Range("A12").Value = 39.78
$ws.Range("B12").Value = 0.16
$ws.Range("C12").Value = 10.07253766059875
$ws.Range("D12").Value = 9.359375

# Row 13
$ws.Range("A13").Value = 39.78
$ws.Range("B13").Value = 0.16
$ws.Range("C13").Value = 12.94005417823792
$ws.Range("D13").Value = 12.25

# Row 14
$ws.Range("A14").Value = 39.06
$ws.Range("B14").Value = 0.04000000000000001
$ws.Range("C14").Value = 3.162457466125488
$ws.Range("D14").Value = 2.515625

# Row 15
$ws.Range("A15").Value = 39.06
$ws.Range("B15").Value = 0.04000000000000001
$ws.Range("C15").Value = 3.130092859268188
$ws.Range("D15").Value = 2.4375

# Row 16
$ws.Range("A16").Value = 39.06
$ws.Range("B16").Value = 0.04000000000000001
$ws.Range("C16").Value = 7.798795700073242
$ws.Range("D16").Value = 6.578125

# Row 17
$ws.Range("A17").Value = 39.06
$ws.Range("B17").Value = 0.04000000000000001
$ws.Range("C17").Value = 3.232581853866577
$ws.Range("D17").Value = 2.546875

# Row 18
$ws.Range("A18").Value = 39.06
$ws.Range("B18").Value = 0.04000000000000001
$ws.Range("C18").Value = 3.882983684539795
$ws.Range("D18").Value = 3.203125

# Row 19
$ws.Range("A19").Value = 39.06
$ws.Range("B19").Value = 0.04000000000000001
$ws.Range("C19").Value = 3.353356122970581
$ws.Range("D19").Value = 2.578125

# Row 20
$ws.Range("A20").Value = 39.06
$ws.Range("B20").Value = 0.04000000000000001
$ws.Range("C20").Value = 3.265630006790161
$ws.Range("D20").Value = 2.625

# Row 21
$ws.Range("A21").Value = 39.06
$ws.Range("B21").Value = 0.04000000000000001
$ws.Range("C21").Value = 3.232125759124756
$ws.Range("D21").Value = 2.4375

# Row 22
$ws.Range("A22").Value = 39.78
$ws.Range("B22").Value = 0.16
$ws.Range("C22").Value = 10.40235781669617
$ws.Range("D22").Value = 9.59375

# Row 23
$ws.Range("A23").Value = 39.78
$ws.Range("B23").Value = 0.16
$ws.Range("C23").Value = 12.70899796485901
$ws.Range("D23").Value = 11.65625

# Row 24
$ws.Range("A24").Value = 39.06
$ws.Range("B24").Value = 0.04000000000000001
$ws.Range("C24").Value = 3.043107748031616
$ws.Range("D24").Value = 2.4375

# Row 25
$ws.Range("A25").Value = 39.78
$ws.Range("B25").Value = 0.16
$ws.Range("C25").Value = 9.73159646987915
$ws.Range("D25").Value = 9.046875

# Row 26
$ws.Range("A26").Value = 39.78
$ws.Range("B26").Value = 0.16
$ws.Range("C26").Value = 14.91842198371887
$ws.Range("D26").Value = 12.90625

# Row 27
$ws.Range("A27").Value = 39.78
$ws.Range("B27").Value = 0.16
$ws.Range("C27").Value = 35.30180263519287
$ws.Range("D27").Value = 19.0625

# Row 28
$ws.Range("A28").Value = 39.78
$ws.Range("B28").Value = 0.16
$ws.Range("C28").Value = 12.96200776100159
$ws.Range("D28").Value = 11.984375

# Row 29
$ws.Range("A29").Value = 39.78
$ws.Range("B29").Value = 0.16
$ws.Range("C29").Value = 10.17955327033997
$ws.Range("D29").Value = 9.296875

# Update selection to A2 as in the target state
$ws.Range("A2").Select()
